$d = $word.ActiveDocument

# --- 1. Move the "_GoBack" bookmark from the middle of the "providing this
#        for patients ... care home." paragraph to the very end of the
#        following "As a group ..." paragraph (just before its paragraph
#        mark), matching the target OOXML. ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

$conclusionPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "As a group we have decided*") {
        $conclusionPara = $p
        break
    }
}
if ($conclusionPara -ne $null) {
    $r = $conclusionPara.Range.Duplicate()
    $r.MoveEnd(1, -1) | Out-Null       # exclude the paragraph mark
    $r.Collapse(0) | Out-Null          # collapse to end (wdCollapseEnd = 0)

    # Placing a bookmark directly on a collapsed range that sits right
    # before a paragraph mark needs a small workaround: insert a
    # throwaway marker character, wrap the bookmark around it, then
    # delete the marker -- the bookmark collapses cleanly in place.
    $r.InsertBefore([char]1)
    $markerRange = $d.Range($r.Start, $r.Start + 1)
    $d.Bookmarks.Add("_GoBack", $markerRange) | Out-Null
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Range.Delete()
}

# --- 2. Delete the long run of trailing blank paragraphs (and the stray
#        tab / page-break paragraphs) at the very end of the document,
#        leaving the section properties intact. ---
$startPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "As a group we have decided*") {
        $startPara = $p.Next()
        break
    }
}
if ($startPara -ne $null) {
    $delRange = $d.Range($startPara.Range.Start, $d.Content.End)
    $delRange.Delete()
}
